$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each target cell holds plain text in the source workbook (prices/
# percentages formatted as strings, e.g. "587.02", "  +0.23%  "). Excel
# COM auto-converts numeric-looking text typed into .Value, so force the
# cell to Text format first, write the value, then restore the original
# (unstyled) look so no stray number-format style gets left behind.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "64.283.76"
Set-TextValue $ws.Range("E2") "  +0.13%  "
Set-TextValue $ws.Range("D3") "3.491.07"
Set-TextValue $ws.Range("E3") "  -1.18%  "
Set-TextValue $ws.Range("E4") "  +0.07%  "
Set-TextValue $ws.Range("D5") "587.30"
Set-TextValue $ws.Range("E5") "  +0.21%  "
Set-TextValue $ws.Range("D6") "134.44"
Set-TextValue $ws.Range("E6") "  +1.05%  "
Set-TextValue $ws.Range("E7") "  -0.01%  "
Set-TextValue $ws.Range("D8") "0.488"
Set-TextValue $ws.Range("E8") "  +0.31%  "
Set-TextValue $ws.Range("E9") "  -0.09%  "
Set-TextValue $ws.Range("D10") "7.25"
Set-TextValue $ws.Range("E10") "  +1.94%  "
Set-TextValue $ws.Range("E11") "  +1.76%  "
Set-TextValue $ws.Range("D12") "4.082.56"
Set-TextValue $ws.Range("E12") "  -1.11%  "
Set-TextValue $ws.Range("E13") "  +1.47%  "
Set-TextValue $ws.Range("E14") "  +0.95%  "
Set-TextValue $ws.Range("D15") "3.490.51"
Set-TextValue $ws.Range("E15") "  -0.92%  "
Set-TextValue $ws.Range("D16") "25.84"
Set-TextValue $ws.Range("E16") "  -6.62%  "
Set-TextValue $ws.Range("D17") "64.338.85"
Set-TextValue $ws.Range("E17") "  +0.30%  "
Set-TextValue $ws.Range("E18") "  +0.33%  "
Set-TextValue $ws.Range("D19") "5.75"
Set-TextValue $ws.Range("E19") "  +2.21%  "
Set-TextValue $ws.Range("E20") "  -3.77%  "
Set-TextValue $ws.Range("D21") "394.49"
Set-TextValue $ws.Range("E21") "  +2.41%  "
Set-TextValue $ws.Range("E22") "  -0.89%  "
Set-TextValue $ws.Range("D23") "3.629.45"
Set-TextValue $ws.Range("E23") "  -1.08%  "
Set-TextValue $ws.Range("D24") "74.73"
Set-TextValue $ws.Range("E24") "  +0.94%  "
Set-TextValue $ws.Range("E25") "  -0.06%  "
Set-TextValue $ws.Range("D26") "5.71"
Set-TextValue $ws.Range("E26") "  +0.85%  "
Set-TextValue $ws.Range("E27") "  -0.73%  "
Set-TextValue $ws.Range("D28") "0.998"
Set-TextValue $ws.Range("E28") "  -0.10%  "
Set-TextValue $ws.Range("E29") "  -1.48%  "
Set-TextValue $ws.Range("B30") "PancakeSwap"
Set-TextValue $ws.Range("C30") "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextValue $ws.Range("D30") "2.24"
Set-TextValue $ws.Range("E30") "  -0.05%  "
Set-TextValue $ws.Range("B31") "Fetch.AI"
Set-TextValue $ws.Range("C31") "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-TextValue $ws.Range("D31") "1.49"
Set-TextValue $ws.Range("E31") "  -5.97%  "
Set-TextValue $ws.Range("D32") "8.25"
Set-TextValue $ws.Range("E32") "  -1.03%  "
Set-TextValue $ws.Range("D33") "3.511.80"
Set-TextValue $ws.Range("E33") "  -0.76%  "
Set-TextValue $ws.Range("E34") "  +3.62%  "
Set-TextValue $ws.Range("E35") "  +0.08%  "
Set-TextValue $ws.Range("D36") "23.43"
Set-TextValue $ws.Range("E36") "  -0.59%  "
Set-TextValue $ws.Range("D37") "5.15"
Set-TextValue $ws.Range("E37") "  -4.31%  "
Set-TextValue $ws.Range("E38") "  -0.14%  "
Set-TextValue $ws.Range("E39") "  -1.49%  "
Set-TextValue $ws.Range("D40") "166.31"
Set-TextValue $ws.Range("E40") "  +4.38%  "
Set-TextValue $ws.Range("D41") "0.0782"
Set-TextValue $ws.Range("E41") "  -1.32%  "
Set-TextValue $ws.Range("D42") "0.806"
Set-TextValue $ws.Range("E42") "  -1.15%  "
Set-TextValue $ws.Range("E43") "  +0.07%  "
Set-TextValue $ws.Range("D44") "25.30"
Set-TextValue $ws.Range("E44") "  -3.23%  "
Set-TextValue $ws.Range("E45") "  -0.63%  "
Set-TextValue $ws.Range("E46") "  +2.09%  "
Set-TextValue $ws.Range("E47") "  -3.69%  "
Set-TextValue $ws.Range("D48") "2.458.13"
Set-TextValue $ws.Range("E48") "  -0.27%  "
Set-TextValue $ws.Range("E49") "  -0.99%  "
Set-TextValue $ws.Range("E50") "  -1.63%  "
Set-TextValue $ws.Range("D51") "0.0261"
Set-TextValue $ws.Range("E51") "  -1.33%  "
